$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.694.87"
Set-TextValue $ws.Range("E2") "  -0.69%  "
Set-TextValue $ws.Range("D3") "3.450.17"
Set-TextValue $ws.Range("E3") "  -2.44%  "
Set-TextValue $ws.Range("D5") "591.10"
Set-TextValue $ws.Range("E5") "  -1.85%  "
Set-TextValue $ws.Range("D6") "178.95"
Set-TextValue $ws.Range("E6") "  -2.47%  "
Set-TextValue $ws.Range("D7") "0.611"
Set-TextValue $ws.Range("E7") "  +2.05%  "
Set-TextValue $ws.Range("E8") "  +0.07%  "
Set-TextValue $ws.Range("D9") "3.447.27"
Set-TextValue $ws.Range("E9") "  -2.62%  "
Set-TextValue $ws.Range("D10") "0.138"
Set-TextValue $ws.Range("E10") "  -1.33%  "
Set-TextValue $ws.Range("D11") "6.95"
Set-TextValue $ws.Range("E11") "  -2.65%  "
Set-TextValue $ws.Range("D12") "0.427"
Set-TextValue $ws.Range("E12") "  -3.20%  "
Set-TextValue $ws.Range("D13") "4.049.97"
Set-TextValue $ws.Range("E13") "  -2.23%  "
Set-TextValue $ws.Range("D14") "32.01"
Set-TextValue $ws.Range("E14") "  -2.44%  "
Set-TextValue $ws.Range("E15") "  -1.16%  "
Set-TextValue $ws.Range("D16") "67.737.07"
Set-TextValue $ws.Range("E16") "  -0.60%  "
Set-TextValue $ws.Range("D17") "0.0000175"
Set-TextValue $ws.Range("E17") "  -4.29%  "
Set-TextValue $ws.Range("D18") "3.454.13"
Set-TextValue $ws.Range("E18") "  -2.06%  "
Set-TextValue $ws.Range("D19") "6.13"
Set-TextValue $ws.Range("E19") "  -4.71%  "
Set-TextValue $ws.Range("D20") "13.97"
Set-TextValue $ws.Range("E20") "  -6.11%  "
Set-TextValue $ws.Range("D21") "391.29"
Set-TextValue $ws.Range("E21") "  -2.13%  "
Set-TextValue $ws.Range("D22") "7.85"
Set-TextValue $ws.Range("E22") "  -3.55%  "
Set-TextValue $ws.Range("E23") "  +2.35%  "
Set-TextValue $ws.Range("D24") "0.999"
Set-TextValue $ws.Range("E24") "  -0.29%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D25") "71.64"
Set-TextValue $ws.Range("E25") "  -2.79%  "
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D26") "0.532"
Set-TextValue $ws.Range("E26") "  -3.00%  "
Set-TextValue $ws.Range("D27") "0.0000119"
Set-TextValue $ws.Range("E27") "  -6.26%  "
Set-TextValue $ws.Range("D28") "10.24"
Set-TextValue $ws.Range("E28") "  -4.61%  "
Set-TextValue $ws.Range("D29") "0.175"
Set-TextValue $ws.Range("E29") "  -2.19%  "
Set-TextValue $ws.Range("E30") "  +0.39%  "
Set-TextValue $ws.Range("D31") "6.05"
Set-TextValue $ws.Range("E31") "  -5.00%  "
Set-TextValue $ws.Range("E32") "  -2.04%  "
Set-TextValue $ws.Range("D33") "1.38"
Set-TextValue $ws.Range("E33") "  -6.83%  "
Set-TextValue $ws.Range("D34") "23.38"
Set-TextValue $ws.Range("E34") "  -3.26%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D35") "1.00"
Set-TextValue $ws.Range("E35") "  -0.05%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D36") "7.24"
Set-TextValue $ws.Range("E36") "  -3.88%  "
Set-TextValue $ws.Range("D37") "1.55"
Set-TextValue $ws.Range("E37") "  -8.40%  "
Set-TextValue $ws.Range("D38") "161.31"
Set-TextValue $ws.Range("E38") "  -1.88%  "
Set-TextValue $ws.Range("D39") "0.884"
Set-TextValue $ws.Range("E39") "  +0.15%  "
Set-TextValue $ws.Range("D40") "2.74"
Set-TextValue $ws.Range("E40") "  -2.69%  "
Set-TextValue $ws.Range("D41") "1.85"
Set-TextValue $ws.Range("E41") "  -6.30%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D42") "6.68"
Set-TextValue $ws.Range("E42") "  -6.79%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D43") "4.60"
Set-TextValue $ws.Range("E43") "  -3.77%  "
Set-TextValue $ws.Range("D44") "25.77"
Set-TextValue $ws.Range("E44") "  -5.30%  "
Set-TextValue $ws.Range("D45") "0.0712"
Set-TextValue $ws.Range("E45") "  -4.53%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D46") "2.704.87"
Set-TextValue $ws.Range("E46") "  -6.13%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D47") "25.88"
Set-TextValue $ws.Range("E47") "  -7.27%  "
Set-TextValue $ws.Range("D48") "41.21"
Set-TextValue $ws.Range("E48") "  -2.89%  "
Set-TextValue $ws.Range("D49") "0.0296"
Set-TextValue $ws.Range("E49") "  -4.11%  "
Set-TextValue $ws.Range("D50") "326.90"
Set-TextValue $ws.Range("E50") "  -7.21%  "
Set-TextValue $ws.Range("D51") "1.03"
Set-TextValue $ws.Range("E51") "  -6.26%  "
